# Automatische test-sync: 2025-08-28 20:11:50
# Append a new "Retour status" log row (row 12) to the Logs sheet and
# bump the Dashboard's "Retour / Terugbetaling" count from 10 to 11.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$row = 12
$logs.Cells.Item($row, 1).Value = "Retour status"
$logs.Cells.Item($row, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($row, 4).Value = "Retour / Terugbetaling"
$logs.Cells.Item($row, 6).Value = "2025-08-28 20:11:08"
$logs.Cells.Item($row, 7).Value = "Ja"
$logs.Cells.Item($row, 8).Value = "Nee"
$logs.Cells.Item($row, 9).Value = "Nee"
$logs.Cells.Item($row, 10).Value = "Nee"

# Update the rolling count on the Dashboard sheet.
$dashboard.Cells.Item(2, 2).Value = 11

# Extend the conditional-formatting ranges so the new row is covered too
# (D2:D11 -> D2:D12, G2:G11 -> G2:G12, H2:H11 -> H2:H12, I2:I11 -> I2:I12,
#  J2:J11 -> J2:J12).
$cfColumns = @("D", "G", "H", "I", "J")
foreach ($col in $cfColumns) {
    $oldRange = $logs.Range("$col" + "2:" + "$col" + "11")
    $newRange = $logs.Range("$col" + "2:" + "$col" + "12")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}
